$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.479.27'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.572.11'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.002'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '291.50'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '49.98'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.84%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3395'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07551'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.61%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.32'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.34%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.046'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.963'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.11%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.570.61'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.36%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.74%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '90.71'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06775'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.295'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.41'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.19'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '22.467.48'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.356'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.624'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.04'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '149.39'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.051'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.31'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.91%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.746.89'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.076'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +7.91%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.252'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.012'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.782'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08361'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.53%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2306'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.340'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.50%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06533'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.460'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.36'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6231'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.44%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.03'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.814'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.93%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.71'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +4.35%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.074'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.219'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.20%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.23%  '
